# "存款" (Deposits) sheet restructuring:
#  - Row 1 becomes a proper header row (field names) instead of a duplicate data row
#  - New columns G:M are added with cash/deposit record metadata for every row
#  - Existing columns B:F keep their original values

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")
$wsLand = $wb.Worksheets.Item("土地")

# ---------------------------------------------------------------------------
# Header row (row 1) - turn into field-name labels
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# Give the new header cells (G1:M1) the same style as the rest of the header row
$ws.Range("B1").Copy()
$ws.Range("G1:M1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Data rows (2-16) - keep B:F as-is, add new columns G:M
# ---------------------------------------------------------------------------
$bank        = @("台中商業銀行","台中商業銀行","台中商業銀行","台中商業銀行","台中商業銀行","台中商業銀行","彰化商業銀行","合作金庫商業銀行","第一商業銀行","臺灣銀行","中華郵政股份有限公司","合作金庫商業銀行","中華郵政股份有限公司","臺灣銀行","台北富邦商業銀行")
$depositType = @("綜合存款","支票存款","活期存款","活期存款","活期儲蓄存款","活期儲蓄存款","活期儲蓄存款","活期儲蓄存款","活期儲蓄存款","綜合存款","活期存款","活期儲蓄存款","綜合存款","綜合存款","綜合存款")
$currency    = @("新臺幣","新臺幣","新臺幣","新臺幣","新臺幣","新臺幣","新臺幣","新臺幣","新臺幣","新臺幣","新臺幣","新臺幣","新臺幣","新臺幣","新臺幣")
$owner       = @("王琴賀","王琴賀","王琴賀","王琴賀","王琴賀","王琴賀","王琴賀","王琴賀","王琴賀","王琴賀","王琴賀","蔡煌瑯","蔡煌瑯","蔡煌瑯","蔡煌瑯")
$total       = @(359122,233640,1099,73427,300,387,23682,9442,61315,42119,374281,60809,672057,713488,283689)
$indexVals   = @(55,56,57,58,59,60,61,62,63,64,65,66,67,68,69)

# A cell elsewhere in the workbook already holds the literal text "2012-04-27"
# as a shared string (土地!K2). Pasting its *value* avoids Excel's automatic
# date-parsing that would otherwise turn a typed "2012-04-27" into a date serial.
$wsLand.Cells.Item(2, 11).Copy()

for ($i = 0; $i -lt 15; $i++) {
    $r = $i + 2

    $ws.Cells.Item($r, 2).Value = $bank[$i]
    $ws.Cells.Item($r, 3).Value = $depositType[$i]
    $ws.Cells.Item($r, 4).Value = $currency[$i]
    $ws.Cells.Item($r, 5).Value = $owner[$i]
    $ws.Cells.Item($r, 6).Value = $total[$i]

    $ws.Cells.Item($r, 7).Value = "deposit"
    $ws.Cells.Item($r, 8).Value = "normal"
    $ws.Cells.Item($r, 9).PasteSpecial(-4163)
    $ws.Cells.Item($r, 10).Value = "蔡煌瑯"
    $ws.Cells.Item($r, 11).Value = 752
    $ws.Cells.Item($r, 12).Value = "tmpd4981"
    $ws.Cells.Item($r, 13).Value = $indexVals[$i]
}

# Give the new data cells (G2:M16) the same style as the rest of the data rows
$ws.Range("B2").Copy()
$ws.Range("G2:M16").PasteSpecial(-4122)

$excel.CutCopyMode = 0
